$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = "'66.253.91"
$ws.Range('D2').Style = 'Normal'
$ws.Range('E2').Value = "'  -0.52%  "
$ws.Range('E2').Style = 'Normal'

$ws.Range('D3').Value = "'3.336.87"
$ws.Range('D3').Style = 'Normal'
$ws.Range('E3').Value = "'  -0.45%  "
$ws.Range('E3').Style = 'Normal'

$ws.Range('D4').Value = "'0.998"
$ws.Range('D4').Style = 'Normal'
$ws.Range('E4').Value = "'  -0.23%  "
$ws.Range('E4').Style = 'Normal'

$ws.Range('D5').Value = "'583.71"
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = "'  +2.90%  "
$ws.Range('E5').Style = 'Normal'

$ws.Range('D6').Value = "'185.66"
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = "'  -2.76%  "
$ws.Range('E6').Style = 'Normal'

$ws.Range('E7').Value = "'  +0.04%  "
$ws.Range('E7').Style = 'Normal'

$ws.Range('D8').Value = "'3.333.86"
$ws.Range('D8').Style = 'Normal'
$ws.Range('E8').Value = "'  -0.20%  "
$ws.Range('E8').Style = 'Normal'

$ws.Range('E9').Value = "'  -2.48%  "
$ws.Range('E9').Style = 'Normal'

$ws.Range('D10').Value = "'0.182"
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = "'  -3.05%  "
$ws.Range('E10').Style = 'Normal'

$ws.Range('D11').Value = "'0.581"
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = "'  -1.85%  "
$ws.Range('E11').Style = 'Normal'

$ws.Range('D12').Value = "'47.20"
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = "'  -1.80%  "
$ws.Range('E12').Style = 'Normal'

$ws.Range('D13').Value = "'0.0000269"
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = "'  -1.84%  "
$ws.Range('E13').Style = 'Normal'

$ws.Range('D14').Value = "'679.78"
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = "'  +11.14%  "
$ws.Range('E14').Style = 'Normal'

$ws.Range('D15').Value = "'3.863.96"
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = "'  -0.52%  "
$ws.Range('E15').Style = 'Normal'

$ws.Range('D16').Value = "'8.51"
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = "'  -2.66%  "
$ws.Range('E16').Style = 'Normal'

$ws.Range('D17').Value = "'66.344.42"
$ws.Range('D17').Style = 'Normal'
$ws.Range('E17').Value = "'  -0.40%  "
$ws.Range('E17').Style = 'Normal'

$ws.Range('D18').Value = "'17.92"
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').Value = "'  -1.53%  "
$ws.Range('E18').Style = 'Normal'

$ws.Range('E19').Value = "'  -0.48%  "
$ws.Range('E19').Style = 'Normal'

$ws.Range('D20').Value = "'3.332.68"
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = "'  -0.61%  "
$ws.Range('E20').Style = 'Normal'

$ws.Range('E21').Value = "'  -1.03%  "
$ws.Range('E21').Style = 'Normal'

$ws.Range('D22').Value = "'0.899"
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = "'  -2.23%  "
$ws.Range('E22').Style = 'Normal'

$ws.Range('D23').Value = "'17.94"
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = "'  -3.77%  "
$ws.Range('E23').Style = 'Normal'

$ws.Range('D24').Value = "'103.23"
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = "'  +1.41%  "
$ws.Range('E24').Style = 'Normal'

$ws.Range('D25').Value = "'5.05"
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = "'  -2.59%  "
$ws.Range('E25').Style = 'Normal'

$ws.Range('D26').Value = "'3.98"
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = "'  -1.89%  "
$ws.Range('E26').Style = 'Normal'

$ws.Range('E27').Value = "'  +0.23%  "
$ws.Range('E27').Style = 'Normal'

$ws.Range('D28').Value = "'9.53"
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = "'  -2.93%  "
$ws.Range('E28').Style = 'Normal'

$ws.Range('D29').Value = "'32.69"
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = "'  +6.48%  "
$ws.Range('E29').Style = 'Normal'

$ws.Range('D30').Value = "'8.52"
$ws.Range('D30').Style = 'Normal'

$ws.Range('E31').Value = "'  -0.78%  "
$ws.Range('E31').Style = 'Normal'

$ws.Range('D32').Value = "'609.23"
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = "'  +7.32%  "
$ws.Range('E32').Style = 'Normal'

$ws.Range('D33').Value = "'3.95"
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = "'  -3.27%  "
$ws.Range('E33').Style = 'Normal'

$ws.Range('D34').Value = "'11.02"
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = "'  -1.38%  "
$ws.Range('E34').Style = 'Normal'

$ws.Range('E35').Value = "'  -1.19%  "
$ws.Range('E35').Style = 'Normal'

$ws.Range('D36').Value = "'3.811.46"
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = "'  +1.58%  "
$ws.Range('E36').Style = 'Normal'

$ws.Range('E37').Value = "'  +0.07%  "
$ws.Range('E37').Style = 'Normal'

$ws.Range('D38').Value = "'56.06"
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = "'  -2.55%  "
$ws.Range('E38').Style = 'Normal'

$ws.Range('D39').Value = "'2.69"
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = "'  -1.92%  "
$ws.Range('E39').Style = 'Normal'

$ws.Range('D40').Value = "'0.0₃0701"
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = "'  -4.77%  "
$ws.Range('E40').Style = 'Normal'

$ws.Range('E41').Value = "'  -3.99%  "
$ws.Range('E41').Style = 'Normal'

$ws.Range('B42').Value = "'Stacks"
$ws.Range('B42').Style = 'Normal'
$ws.Range('C42').Value = "'https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range('C42').Style = 'Normal'
$ws.Range('D42').Value = "'3.20"
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = "'  -4.14%  "
$ws.Range('E42').Style = 'Normal'

$ws.Range('B43').Value = "'InjectiveProtocol"
$ws.Range('B43').Style = 'Normal'
$ws.Range('C43').Value = "'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range('C43').Style = 'Normal'
$ws.Range('D43').Value = "'32.83"
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = "'  -4.41%  "
$ws.Range('E43').Style = 'Normal'

$ws.Range('E44').Value = "'  +5.16%  "
$ws.Range('E44').Style = 'Normal'

$ws.Range('E45').Value = "'  -2.66%  "
$ws.Range('E45').Style = 'Normal'

$ws.Range('D46').Value = "'0.0416"
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = "'  -3.06%  "
$ws.Range('E46').Style = 'Normal'

$ws.Range('D47').Value = "'3.01"
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = "'  -11.75%  "
$ws.Range('E47').Style = 'Normal'

$ws.Range('E48').Value = "'  -2.03%  "
$ws.Range('E48').Style = 'Normal'

$ws.Range('B49').Value = "'ThetaToken"
$ws.Range('B49').Style = 'Normal'
$ws.Range('C49').Value = "'https://coinranking.com/coin/B42IRxNtoYmwK+thetatoken-theta"
$ws.Range('C49').Style = 'Normal'
$ws.Range('D49').Value = "'2.57"
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = "'  -2.23%  "
$ws.Range('E49').Style = 'Normal'

$ws.Range('B50').Value = "'FirstDigitalUSD"
$ws.Range('B50').Style = 'Normal'
$ws.Range('C50').Value = "'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
$ws.Range('C50').Style = 'Normal'
$ws.Range('D50').Value = "'1.00"
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = "'  +0.31%  "
$ws.Range('E50').Style = 'Normal'

$ws.Range('B51').Value = "'Monero"
$ws.Range('B51').Style = 'Normal'
$ws.Range('C51').Value = "'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range('C51').Style = 'Normal'
$ws.Range('D51').Value = "'130.84"
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = "'  +5.54%  "
$ws.Range('E51').Style = 'Normal'

Write-Host "Applied all cryptos updates"